# Apply the STGraph "Overview" slide edits:
#  - paragraph "Usually built on index-free adjacency principles." gains a
#    bold run around "index-free adjacency principles" and extra lead-in text.
#  - paragraph "Usually built on Log-Structured Merge Trees." gains a bold
#    run around "Log-Structured Merge Tree" and a new trailing clause.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- "index-free adjacency principles" bullet (2nd paragraph) ---
$para1 = $tr.Paragraphs(2, 1)
$run1 = $para1.Runs(1, 1)
$run1.Text = "Usually built on data structures implementing index-free adjacency principles."
$bold1 = $run1.Characters(47, 31)
$bold1.Font.Bold = $true

# --- "Log-Structured Merge Tree" bullet (4th paragraph) ---
$para2 = $tr.Paragraphs(4, 1)
$run2 = $para2.Runs(1, 1)
$run2.Text = "Usually built on Log-Structured Merge Tree based data structures."
$bold2 = $run2.Characters(18, 25)
$bold2.Font.Bold = $true
